$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(9)
$r = $p.Range
$r.Collapse(0)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="75D7672F" w14:textId="1D721412" w:rsidR="00FD5269" w:rsidRDefault="00FD5269" w:rsidP="00281FB5"><w:r><w:t xml:space="preserve">(We can get </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00C84DD5"><w:t>R</w:t></w:r><w:r><w:t>_t</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> values per state)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>Psuedo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t>-code</w:t></w:r></w:p><w:p><w:r><w:t>Extract data from rt.csv</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Each state </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>R_t</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> value starts from a different date and goes to the day before the current date. Have to find a way to reconcile that. If we do that then the process is simple</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Anyways extract average of each state&#8217;s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>R_t</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from each day, and average that average to create a US-wide average. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Plug into model and profit. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r.InsertXML($xml)
Write-Output "done"
